$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Added PDQ page load tests:
#   - new "PDQPage" worksheet (page-load rows for the new PDQ Cancer Info
#     Summary content type), inserted right before "TopicPage"
#   - new blank "Sheet2" worksheet, inserted right after "PDQPage"
# ---------------------------------------------------------------------------

$topicPage = $wb.Worksheets.Item("TopicPage")

# Duplicate TopicPage (Before:=TopicPage) so the new sheet starts out with the
# same header style/column setup, then rename it and swap in the PDQ data.
$topicPage.Copy($topicPage, $null) | Out-Null
$pdqPage = $wb.Worksheets.Item("TopicPage (2)")
$pdqPage.Name = "PDQPage"

# New blank worksheet, right after PDQPage (and before TopicPage).
$sheet2 = $wb.Worksheets.Add($null, $pdqPage)
$sheet2.Name = "Sheet2"

# Clear the body rows inherited from TopicPage (keep header row 1 as-is).
$pdqPage.Rows.Item(2).Resize(4, 1).EntireRow.ClearContents() | Out-Null

$pdqRows = @(
  @("/types/lung/hp/lung-prevention-pdq", "PDQ Cancer Info Summary"),
  @("/espanol/tipos/pulmon/pro/prevencion-pulmon-pdq", "PDQ Cancer Info Summary"),
  @("/types/lung/hp/lung-prevention-pdq#section/all", "PDQ Cancer Info Summary"),
  @("/types/lung/hp/lung-prevention-pdq#section/_16", "PDQ Cancer Info Summary"),
  @("/types/lung/hp/lung-prevention-pdq#link/_225_toc", "PDQ Cancer Info Summary"),
  @("/espanol/cancer/deteccion/aspectos-generales-deteccion-paciente-pdq", "PDQ Cancer Info Summary"),
  @("/about-cancer/screening/patient-screening-overview-pdq", "PDQ Cancer Info Summary"),
  @("/espanol/cancer/deteccion/aspectos-generales-deteccion-paciente-pdq#section/all", "PDQ Cancer Info Summary"),
  @("/espanol/cancer/deteccion/aspectos-generales-deteccion-paciente-pdq#section/_149", "PDQ Cancer Info Summary"),
  @("/espanol/cancer/deteccion/aspectos-generales-deteccion-paciente-pdq#link/_8", "PDQ Cancer Info Summary")
)

$r = 2
foreach ($row in $pdqRows) {
  $pdqPage.Cells.Item($r, 1).Value = $row[0]
  $pdqPage.Cells.Item($r, 2).Value = $row[1]
  $r = $r + 1
}

# Resize the columns to fit the (now much longer) PDQ paths/labels.
$pdqPage.Columns.Item(1).AutoFit() | Out-Null
$pdqPage.Columns.Item(2).AutoFit() | Out-Null

# PDQPage ends up the active/selected tab, cursor parked just below the data.
$pdqPage.Activate()
$pdqPage.Range("A12").Select() | Out-Null
